# Updated cryptos list on Sat Sep 28 03:56:57 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns on Sheet1 with new
# quotes, and swaps the Bittensor / Binance-PegBSC-USD rows (31/32).
#
# Several Price values look like plain numbers to Excel's auto-detection
# (e.g. "612.94", "1.00"), which would silently coerce them to numeric
# cells on assignment. Those are written with a leading apostrophe to
# force text, then the resulting "quote prefix" style is cleared via
# Style = "Normal" so the cell keeps its original (default) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.081.56'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '2.691.74'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'612.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = "'158.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("E9").Value = '  +6.28%  '
$ws.Range("D10").Value = "'6.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.40%  '
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("E13").Value = '  +10.43%  '
$ws.Range("D14").Value = "'30.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.78%  '
$ws.Range("D15").Value = '3.173.34'
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("D16").Value = '65.934.19'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = '2.690.21'
$ws.Range("E17").Value = '  +1.78%  '
$ws.Range("D18").Value = "'12.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.96%  '
$ws.Range("D19").Value = "'4.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = "'7.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.87%  '
$ws.Range("D21").Value = "'359.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").Value = "'71.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.58%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("E24").Value = '  +16.94%  '
$ws.Range("E25").Value = '  +5.13%  '
$ws.Range("E26").Value = '  -1.89%  '
$ws.Range("E27").Value = '  +0.54%  '
$ws.Range("E28").Value = '  +3.99%  '
$ws.Range("D29").Value = "'8.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.55%  '
$ws.Range("E30").Value = '  +1.80%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = "'540.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.34%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("E34").Value = '  +4.82%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("D37").Value = "'20.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("D38").Value = "'164.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("D39").Value = "'1.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = "'168.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.78%  '
$ws.Range("D43").Value = "'42.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = "'0.0633"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.91%  '
$ws.Range("D46").Value = "'23.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.05%  '
$ws.Range("E47").Value = '  +4.03%  '
$ws.Range("D48").Value = "'0.0268"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.81%  '
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").Value = "'20.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.28%  '
$ws.Range("E51").Value = '  +1.17%  '
